# Automatische test-sync: 2025-06-24 20:04:50
#
# Adds the new "Kan mijn wachtwoord niet resetten" ticket to the Logs
# sheet (row 15) and refreshes the Dashboard summary sheet / chart data
# so the "IT / Technisch probleem" bucket reflects the new count and
# moves into its sorted position.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs sheet: append new row 15 -----------------------------------
$row = 15

$logs.Cells.Item($row, 1).Value = "Kan mijn wachtwoord niet resetten"
$logs.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($row, 3).Value = "Ik krijg geen e-mail bij wachtwoord resetten."
$logs.Cells.Item($row, 4).Value = "IT / Technisch probleem"
$logs.Cells.Item($row, 5).Value = "Beste klant,`nBedankt voor het melden van dit probleem. Om u beter te kunnen helpen, heb ik wat meer informatie van u nodig. Kunt u mij uw gebruikersnaam of het e-mailadres waarmee u probeert uw wachtwoord te resetten doorgeven? Op die manier kunnen we verder onderzoeken wat er misgaat en u zo snel mogelijk helpen.`nMet vriendelijke groet,`n[Naam bedrijf] E-mailassistent"
$logs.Cells.Item($row, 6).Value = "2025-06-24 20:04:03"
$logs.Cells.Item($row, 7).Value = "Ja"

# --- Logs sheet: extend conditional formatting to the new row --------
# Grow the existing rules' AppliesTo range in place (rather than deleting
# and re-adding) so the dxf/style ids and rule ordering are preserved.
$dRules = $logs.Range("D2:D14").FormatConditions
$dRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D15"))

$gRules = $logs.Range("G2:G14").FormatConditions
$gRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G15"))

# --- Dashboard sheet: refresh the category summary --------------------
# New order/counts after the extra "IT / Technisch probleem" ticket:
$dashboard.Cells.Item(5, 1).Value = "IT / Technisch probleem"
$dashboard.Cells.Item(5, 2).Value = 2
$dashboard.Cells.Item(6, 1).Value = "Offerte / Prijsaanvraag"
$dashboard.Cells.Item(6, 2).Value = 1
$dashboard.Cells.Item(7, 1).Value = "Productinformatie"
$dashboard.Cells.Item(7, 2).Value = 1
$dashboard.Cells.Item(8, 1).Value = "Bestelling / Levering"
$dashboard.Cells.Item(8, 2).Value = 1
